# Scheduled-runner market-price refresh for the Siren Profits workbook.
# Each per-job worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) has a Leve-profit
# table (Table_<JOB>) whose H:N columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) are re-pulled from the market-board source on each run. This
# applies the refreshed values row by row per sheet.

$wb = $excel.ActiveWorkbook

# One entry per updated Leve row: Sheet, Row, Updates (column -> new value),
# Clears (columns whose value was removed entirely by this refresh).
$updates = @(
    @{ Sheet = "ALC"; Row = 15; Updates = @{ "H" = 1123.3898; "I" = 1123.3898; "K" = 3370.1694; "M" = -3201.1694 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 33; Updates = @{ "H" = 251.03847; "I" = 255.70833; "K" = 255.70833; "M" = -26.70832999999999 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 38; Updates = @{ "H" = 2357.1667; "J" = 12500; "L" = 37500; "N" = -38244 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 58; Updates = @{ "H" = 3019.476; "J" = 4051.1428; "L" = 12153.4284; "N" = -12453.4284 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 112; Updates = @{ "H" = 30929.5; "J" = 31610.707; "L" = 94832.121; "N" = -97048.121 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 113; Updates = @{ "H" = 11235.75; "J" = 7356.4443; "L" = 7356.4443; "N" = -13864.4443 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 137; Updates = @{ "H" = 652162.8; "I" = 1473743.6; "J" = 13155.556; "K" = 4421230.800000001; "L" = 39466.66800000001; "M" = -4418680.800000001; "N" = -44566.66800000001 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 138; Updates = @{ "H" = 4840.275; "J" = 5322.282; "L" = 15966.846; "N" = -26246.846 }; Clears = @() },
    @{ Sheet = "ALC"; Row = 141; Updates = @{ "H" = 3953.4; "I" = 1878.3572; "K" = 5635.071599999999; "M" = -455.0715999999993 }; Clears = @() },
    @{ Sheet = "ARM"; Row = 32; Updates = @{ "H" = 4512.971; "I" = 4707.758; "K" = 4707.758; "M" = -4420.758 }; Clears = @() },
    @{ Sheet = "ARM"; Row = 45; Updates = @{ "H" = 82532.22; "I" = 308990.44; "K" = 308990.44; "M" = -308613.44 }; Clears = @() },
    @{ Sheet = "ARM"; Row = 80; Updates = @{ "H" = 79700; "J" = 79700; "L" = 79700; "N" = -81696 }; Clears = @() },
    @{ Sheet = "ARM"; Row = 83; Updates = @{ "H" = 79700; "J" = 79700; "L" = 239100; "N" = -249084 }; Clears = @() },
    @{ Sheet = "BSM"; Row = 50; Updates = @{ "H" = 104780; "J" = 104780; "L" = 104780; "N" = -105928 }; Clears = @() },
    @{ Sheet = "BSM"; Row = 52; Updates = @{ "H" = 40999.332; "I" = 30000; "J" = 43199.2; "K" = 30000; "L" = 43199.2; "M" = -29737; "N" = -43725.2 }; Clears = @() },
    @{ Sheet = "BSM"; Row = 55; Updates = @{ "H" = 0; "J" = 0; "L" = 0 }; Clears = @("N") },
    @{ Sheet = "BSM"; Row = 121; Updates = @{ "H" = 40999.332; "I" = 30000; "J" = 43199.2; "K" = 30000; "L" = 43199.2; "M" = -28253; "N" = -46693.2 }; Clears = @() },
    @{ Sheet = "CRP"; Row = 22; Updates = @{ "H" = 456.1; "I" = 355.69232; "K" = 355.69232; "M" = -5.692319999999995 }; Clears = @() },
    @{ Sheet = "CRP"; Row = 86; Updates = @{ "H" = 11847.889; "I" = 11387; "K" = 11387; "M" = -10264 }; Clears = @() },
    @{ Sheet = "CRP"; Row = 89; Updates = @{ "H" = 11847.889; "I" = 11387; "K" = 56935; "M" = -51319 }; Clears = @() },
    @{ Sheet = "CRP"; Row = 114; Updates = @{ "H" = 46237.25; "J" = 54983; "L" = 54983; "N" = -63661 }; Clears = @() },
    @{ Sheet = "CRP"; Row = 117; Updates = @{ "H" = 0; "J" = 0; "L" = 0 }; Clears = @("N") },
    @{ Sheet = "CUL"; Row = 5; Updates = @{ "H" = 358610.72; "J" = 668239.4399999999; "L" = 2004718.32; "N" = -2004942.32 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 34; Updates = @{ "H" = 2167653.2; "J" = 990000; "L" = 2970000; "N" = -2970168 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 40; Updates = @{ "H" = 59.23077; "I" = 41.88889; "J" = 64.433334; "K" = 167.55556; "L" = 257.733336; "M" = -98.55556000000001; "N" = -395.733336 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 47; Updates = @{ "H" = 311.75; "I" = 82.666664; "J" = 999; "K" = 247.999992; "L" = 2997; "M" = 183.000008; "N" = -3859 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 68; Updates = @{ "H" = 14709835; "J" = 19235138; "L" = 57705414; "N" = -57707036 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 71; Updates = @{ "H" = 14709835; "J" = 19235138; "L" = 173116242; "N" = -173124354 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 112; Updates = @{ "H" = 5844.4287; "I" = 5856.636; "J" = 5799.6665; "K" = 17569.908; "L" = 17398.9995; "M" = -16461.908; "N" = -19614.9995 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 113; Updates = @{ "H" = 1491.091; "I" = 407.1; "J" = 2394.4167; "K" = 1221.3; "L" = 7183.250100000001; "M" = 948.6999999999998; "N" = -11523.2501 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 127; Updates = @{ "H" = 1086.5555; "I" = 794; "J" = 1123.125; "K" = 2382; "L" = 3369.375; "M" = 2578; "N" = -13289.375 }; Clears = @() },
    @{ Sheet = "CUL"; Row = 135; Updates = @{ "H" = 358610.72; "J" = 668239.4399999999; "L" = 6014154.959999999; "N" = -6019224.959999999 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 24; Updates = @{ "H" = 733824.75 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 80; Updates = @{ "H" = 12937.6; "J" = 3849; "L" = 3849; "N" = -5845 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 83; Updates = @{ "H" = 12937.6; "J" = 3849; "L" = 19245; "N" = -29229 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 128; Updates = @{ "H" = 79985.5; "J" = 79985.5; "L" = 79985.5; "N" = -89945.5 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 129; Updates = @{ "H" = 42487.5; "J" = 42487.5; "L" = 42487.5; "N" = -52487.5 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 132; Updates = @{ "H" = 6240.8423; "I" = 4622.7646; "K" = 13868.2938; "M" = -11338.2938 }; Clears = @() },
    @{ Sheet = "GSM"; Row = 135; Updates = @{ "H" = 95911; "J" = 95911; "L" = 95911; "N" = -106051 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 40; Updates = @{ "H" = 57291.855; "I" = 90049.625; "J" = 13614.833; "K" = 90049.625; "L" = 13614.833; "M" = -89913.625; "N" = -13886.833 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 68; Updates = @{ "H" = 4799.5; "I" = 2360.2; "J" = 6541.857; "K" = 2360.2; "L" = 6541.857; "M" = -1611.2; "N" = -8039.857 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 71; Updates = @{ "H" = 4799.5; "I" = 2360.2; "J" = 6541.857; "K" = 11801; "L" = 32709.285; "M" = -8057; "N" = -40197.285 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 127; Updates = @{ "H" = 250200430; "J" = 267238.34; "L" = 267238.34; "N" = -277158.34 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 132; Updates = @{ "H" = 1642477.9; "I" = 2110614.5; "K" = 6331843.5; "M" = -6329313.5 }; Clears = @() },
    @{ Sheet = "LTW"; Row = 136; Updates = @{ "H" = 6234.1333; "I" = 5138.25; "J" = 6632.636; "K" = 15414.75; "L" = 19897.908; "M" = -12864.75; "N" = -24997.908 }; Clears = @() },
    @{ Sheet = "WVR"; Row = 20; Updates = @{ "H" = 20252; "I" = 25336; "J" = 5000; "K" = 25336; "L" = 5000; "M" = -25096; "N" = -5480 }; Clears = @() },
    @{ Sheet = "WVR"; Row = 39; Updates = @{ "H" = 30044; "I" = 30044; "K" = 30044; "M" = -29631 }; Clears = @() },
    @{ Sheet = "WVR"; Row = 123; Updates = @{ "H" = 49000; "J" = 49000; "L" = 49000; "N" = -58800 }; Clears = @() },
    @{ Sheet = "WVR"; Row = 132; Updates = @{ "H" = 44497.8; "I" = 55284; "J" = 19330; "K" = 165852; "L" = 57990; "M" = -163322; "N" = -63050 }; Clears = @() },
    @{ Sheet = "WVR"; Row = 136; Updates = @{ "H" = 2312.8333; "I" = 2137.8; "J" = 2604.5557; "K" = 6413.400000000001; "L" = 7813.6671; "M" = -3863.400000000001; "N" = -12913.6671 }; Clears = @() }
)

foreach ($entry in $updates) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($col in $entry.Updates.Keys) {
        $ws.Range("$col$($entry.Row)").Value = $entry.Updates[$col]
    }
    foreach ($col in $entry.Clears) {
        $ws.Range("$col$($entry.Row)").ClearContents()
    }
}
